$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing grade value (2) for the "Render meshes" row in column C
$ws.Range("C20").Value = 2

# Update the active cell selection to C20, matching the saved view state
$ws.Range("C20").Select()
